# Fixed update to excel issue
$wb = $excel.ActiveWorkbook

# Rename header cells on existing sheets
$wsWeekly = $wb.Worksheets.Item("Weekly Quantity")
$wsWeekly.Range("B1").Value = "Weekly_PO_Qty"

$wsMonthly = $wb.Worksheets.Item("Monthly Trend")
$wsMonthly.Range("B1").Value = "Monthly_PO_Qty"

# Add a new sheet "PO Forecast" after "Monthly Trend"
$wsForecast = $wb.Worksheets.Add($null, $wsMonthly)
$wsForecast.Name = "PO Forecast"

# Headers
$wsForecast.Range("A1").Value = "ds"
$wsForecast.Range("B1").Value = "PO_Forecast"
$wsForecast.Range("C1").Value = "yhat_lower"
$wsForecast.Range("D1").Value = "yhat_upper"

# Data rows
$data = @(
    @(44983.99999999999, 8,  -13.54033659815663, 28.27112205901213),
    @(45004.99999999999, 14, -5.47679374089636,  34.25169601461835),
    @(45018.99999999999, 18, -4.10786475883537,  37.97405307847627),
    @(45046.99999999999, 26, 4.656799205620709,  49.04217157487792),
    @(45088.99999999999, 38, 15.50295391447736,  59.42702658261985),
    @(45095.99999999999, 40, 16.53107948053091,  61.00572623170741),
    @(45102.99999999999, 42, 20.33057259120463,  62.87760092590094),
    @(45109.99999999999, 44, 22.10682576205002,  65.84110845020535),
    @(45116.99999999999, 46, 22.92464632241801,  68.44956642931274),
    @(45123.99999999999, 48, 26.03813005306147,  70.37062452389659),
    @(45130.99999999999, 50, 27.92171329825141,  69.87441517712246),
    @(45137.99999999999, 52, 30.01925422432577,  73.04287936733478),
    @(45144.99999999999, 54, 30.49111758512256,  75.4409022868758),
    @(45151.99999999999, 56, 34.46584562674527,  76.39314205672244),
    @(45158.99999999999, 58, 35.644702170016,    79.12049177483895),
    @(45165.99999999999, 60, 36.80869540763415,  80.94878391936219),
    @(45172.99999999999, 62, 39.56681526909038,  83.25561374459986),
    @(45179.99999999999, 64, 40.55544836503702,  85.59315223740235)
)

$r = 2
foreach ($row in $data) {
    $wsForecast.Cells.Item($r, 1).Value = $row[0]
    $wsForecast.Cells.Item($r, 2).Value = $row[1]
    $wsForecast.Cells.Item($r, 3).Value = $row[2]
    $wsForecast.Cells.Item($r, 4).Value = $row[3]
    $r++
}

# Match header/date styles used on the other sheets (bold header row, date-formatted column A)
$header = $wsForecast.Range("A1:D1")
$header.Font.Bold = $true
$header.HorizontalAlignment = -4108
$header.VerticalAlignment = -4160
$header.Borders.LineStyle = 1

$wsForecast.Range("A2:A19").NumberFormat = "YYYY-MM-DD HH:MM:SS"

# Restore the originally active sheet/selection
$wsWeekly.Activate() | Out-Null
$wsWeekly.Range("A1").Select() | Out-Null
